$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value. Values are written as plain text (matching the
# source xlsx, which stores every cell as inlineStr) by forcing a text
# number format before the write, then clearing the format again so the
# cell's style index is left exactly as it was (style 0/general).
$updates = [ordered]@{
    "D2" = "67.008.66"
    "E2" = "  -1.26%  "
    "D3" = "2.609.88"
    "E3" = "  -0.80%  "
    "D4" = "1.00"
    "E4" = "  -0.05%  "
    "D5" = "592.92"
    "E5" = "  -0.94%  "
    "D6" = "152.45"
    "E6" = "  -1.94%  "
    "E7" = "  +0.03%  "
    "D8" = "0.554"
    "E8" = "  +1.69%  "
    "D9" = "2.608.86"
    "E9" = "  -0.78%  "
    "E10" = "  -3.81%  "
    "E11" = "  +0.31%  "
    "D12" = "5.17"
    "E12" = "  -2.33%  "
    "D13" = "0.346"
    "E13" = "  -3.58%  "
    "D14" = "27.50"
    "E14" = "  -1.68%  "
    "D15" = "3.081.66"
    "E15" = "  -0.84%  "
    "E16" = "  -4.87%  "
    "D17" = "66.903.98"
    "E17" = "  -1.21%  "
    "D18" = "2.622.61"
    "E18" = "  -0.31%  "
    "D19" = "365.36"
    "E19" = "  -0.94%  "
    "D20" = "11.02"
    "E20" = "  -3.38%  "
    "D21" = "7.35"
    "E21" = "  -5.69%  "
    "E22" = "  -0.56%  "
    "D23" = "2.05"
    "E23" = "  +0.08%  "
    "E24" = "  +0.14%  "
    "D25" = "10.08"
    "E25" = "  +1.10%  "
    "D26" = "66.17"
    "E26" = "  -1.58%  "
    "D27" = "2.739.43"
    "E27" = "  -1.18%  "
    "D28" = "585.34"
    "E28" = "  +0.62%  "
    "E29" = "  -0.16%  "
    "E30" = "  -4.89%  "
    "E31" = "  -5.21%  "
    "D32" = "7.70"
    "E32" = "  -3.55%  "
    "E33" = "  -3.08%  "
    "E34" = "  +0.05%  "
    "D35" = "0.123"
    "E35" = "  -7.96%  "
    "E36" = "  -3.15%  "
    "D37" = "4.83"
    "E37" = "  -2.66%  "
    "D38" = "154.82"
    "E38" = "  -0.26%  "
    "D39" = "18.89"
    "E39" = "  -3.66%  "
    "D40" = "0.365"
    "E40" = "  -1.72%  "
    "D41" = "5.21"
    "E41" = "  -3.61%  "
    "D42" = "1.80"
    "E42" = "  -4.28%  "
    "D43" = "2.57"
    "E43" = "  -0.71%  "
    "B44" = "WhiteBITCoin"
    "C44" = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
    "D44" = "16.63"
    "E44" = "  +0.62%  "
    "B45" = "OKB"
    "C45" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D45" = "40.76"
    "E45" = "  -1.90%  "
    "E46" = "  -0.05%  "
    "D47" = "154.49"
    "E47" = "  -2.21%  "
    "D48" = "0.0₆0292"
    "E48" = "  -1.10%  "
    "D49" = "3.71"
    "E49" = "  -1.02%  "
    "D50" = "21.48"
    "E50" = "  +2.07%  "
    "D51" = "0.616"
    "E51" = "  -3.32%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.ClearFormats()
}
